$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.356.44"
$ws.Range("E2").Value = "  +4.45%  "

$ws.Range("D3").Value = "1.803.62"
$ws.Range("E3").Value = "  +2.08%  "

$ws.Range("D4").Value = "'0.9997"
$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "'315.81"
$ws.Range("E5").Value = "  +0.99%  "

$ws.Range("D6").Value = "'0.9992"
$ws.Range("E6").Value = "  +0.12%  "

$ws.Range("D7").Value = "'0.5507"
$ws.Range("E7").Value = "  +5.59%  "

$ws.Range("D8").Value = "'0.3858"
$ws.Range("E8").Value = "  +7.05%  "

$ws.Range("D9").Value = "'0.07596"
$ws.Range("E9").Value = "  +3.67%  "

$ws.Range("D10").Value = "'42.63"
$ws.Range("E10").Value = "  +0.70%  "

$ws.Range("E11").Value = "  +3.78%  "

$ws.Range("D12").Value = "'0.9996"
$ws.Range("E12").Value = "  +0.17%  "

$ws.Range("D13").Value = "'21.16"
$ws.Range("E13").Value = "  +3.14%  "

$ws.Range("E14").Value = "  +2.23%  "

$ws.Range("D15").Value = "'7.354"
$ws.Range("E15").Value = "  +5.89%  "

$ws.Range("D16").Value = "1.801.44"
$ws.Range("E16").Value = "  +2.03%  "

$ws.Range("D17").Value = "'92.18"
$ws.Range("E17").Value = "  +4.52%  "

$ws.Range("E18").Value = "  +2.57%  "

$ws.Range("D19").Value = "'0.06435"
$ws.Range("E19").Value = "  +0.36%  "

$ws.Range("E20").Value = "  +0.15%  "

$ws.Range("E21").Value = "  +3.76%  "

$ws.Range("D22").Value = "'5.991"
$ws.Range("E22").Value = "  +2.92%  "

$ws.Range("D23").Value = "28.373.66"
$ws.Range("E23").Value = "  +4.13%  "

$ws.Range("D24").Value = "'11.46"
$ws.Range("E24").Value = "  +1.12%  "

$ws.Range("D25").Value = "'2.133"
$ws.Range("E25").Value = "  +3.32%  "

$ws.Range("D26").Value = "'158.00"
$ws.Range("E26").Value = "  +2.37%  "

$ws.Range("D27").Value = "'20.67"
$ws.Range("E27").Value = "  +3.09%  "

$ws.Range("D28").Value = "'2.388"
$ws.Range("E28").Value = "  +3.65%  "

$ws.Range("D29").Value = "2.009.85"
$ws.Range("E29").Value = "  +2.37%  "

$ws.Range("D30").Value = "'123.82"
$ws.Range("E30").Value = "  +2.49%  "

$ws.Range("D31").Value = "'1.124"
$ws.Range("E31").Value = "  +6.40%  "

$ws.Range("D32").Value = "'0.1020"
$ws.Range("E32").Value = "  +4.88%  "

$ws.Range("D33").Value = "'5.745"
$ws.Range("E33").Value = "  +4.62%  "

$ws.Range("D34").Value = "'3.669"
$ws.Range("E34").Value = "  +2.06%  "

$ws.Range("D35").Value = "'0.2334"
$ws.Range("E35").Value = "  +15.23%  "

$ws.Range("D36").Value = "'0.06340"
$ws.Range("E36").Value = "  +6.08%  "

$ws.Range("D37").Value = "'0.02319"
$ws.Range("E37").Value = "  +4.66%  "

$ws.Range("E38").Value = "  +11.85%  "

$ws.Range("D39").Value = "'11.62"
$ws.Range("E39").Value = "  +3.99%  "

$ws.Range("E40").Value = "  +4.94%  "

$ws.Range("D41").Value = "'0.6411"
$ws.Range("E41").Value = "  +4.94%  "

$ws.Range("D42").Value = "'0.9989"
$ws.Range("E42").Value = "  +0.18%  "

$ws.Range("D43").Value = "'1.155"
$ws.Range("E43").Value = "  +2.70%  "

$ws.Range("D44").Value = "'1.382"
$ws.Range("E44").Value = "  -3.53%  "

$ws.Range("D45").Value = "'13.51"
$ws.Range("E45").Value = "  +2.39%  "

$ws.Range("D46").Value = "'0.5983"
$ws.Range("E46").Value = "  +4.44%  "

$ws.Range("E47").Value = "  +1.86%  "

$ws.Range("D48").Value = "'124.30"
$ws.Range("E48").Value = "  +2.60%  "

$ws.Range("D49").Value = "'1.987"
$ws.Range("E49").Value = "  +5.98%  "

$ws.Range("D50").Value = "'1.149"
$ws.Range("E50").Value = "  +3.76%  "

$ws.Range("D51").Value = "'0.06908"
$ws.Range("E51").Value = "  +3.23%  "
